# Apply Week 15 logging + Week 16 simulation updates to 2021 Team Data workbook
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append new per-game yardage figures ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = '1 11 5 15 3 3 6 2 8 2 10 5 3 9 5 9 3 11 1 1 6 6 14 2 -1 3 2 1 1 1 2 7 4 0 2 6 3 2 2 0 0 4 2 12 0 1 7 1 5 9 -7 8 6 1 0 5 0 7 5 3 6 5 3 3 3 3 8 5 2 7 9 5 7 5 -1 1 -2 1 -3 3 11 3 14 6 3 2 5 5 4 6 10 6 3 3 3 2 5 4 2 3 8 1 3 2 5 8 1 8 5 3 4 6 9 2 9 4 4 4 8 7 1 0 3 2 8 1 5 2 1 5 3 1 23 2 10 5 3 12 3 3 0 11 2 1 8 3 -2 1 5 4 2 3 1 3 4 3 2 0 0 2 18 3 -1 6 2 5 -3 3 9 1 7 1 12 4 -2 20 10 1 4 3 13 3 1 3 -1 3 3 11 7 1 1 3 2 15 2 4 2 3 6 5 1 22 1 3 9 2 9 11 9 9 14 2 4 1 3 8 -1 3 1 6 -1 3 2 44 2 4 2 2 3 1 2 5 10 3 1 6 8 0 2 2 6 4 4 6 13 3 0 5 3 2 6 -1 4 2 14 9 1 -2 4 3 1 9 1 -1 28 3 1 17 7 3 1 2 -2 3 4 8 2 0 0 2 2 1 4 -5 4 2 5 3 -1 3 2 9 9 3 -1 25 24 7 0 4 8 6 5 9 23 4 6 8 0 7 7 5 11 7 1 1 3 10 -1 2 4 7 1 3 6 0 8 9 2 16 5 1 0 -2 7 13 -2 2 2 -1 9 3 0 9 3 8 3 1 -1 2 -1 0 0 44 8 -1 4 -3 5 -2 1 1 9 2 2 6 6 3 -1 3 8 0 0 7 3 -1 4 1 1 -1 3 -3 1 0 -5'
$wsYDS.Range("C2").Value = '-1 1 6 6 3 3 5 3 3 0 3 5 3 2 6 1 1 0 2 5 -1 -1 16 6 0 7 3 5 9 0 2 -2 0 11 0 -1 1 3 2 -5 3 7 -2 4 3 6 1 2 0 1 17 17 2 2 2 0 -1 7 12 0 4 0 1 7 8 1 17 20 1 1 -3 2 4 1 4 2 1 6 6 5 8 3 0 1 8 5 1 -7 13 5 3 3 1 4 1 2 3 3 2 5 6 16 1 9 2 11 2 -2 -2 16 2 1 3 6 0 0 1 8 12 3 -3 1 1 5 12 2 14 3 6 1 4 5 0 1 -1 5 4 -1 -1 -1 4 13 9 -1 2 -1 0 24 6 2 5 6 3 1 0 4 1 5 0 2 2 1 3 4 4 5 1 2 1 0 1 3 -4 -1 2 1 2 4 2 5 4 -1 6 5 7 1 2 1 5 3 2 2 1 1 0 2 3 4 5 6 -5 0 1 2 3 4 -3 4 1 25 4 4 20 4 3 3 2 4 4 18 9 4 13 4 4 4 5 6 0 2 5 6 0 4 -3 7 5 10 6 2 -6 3 -1 7 2 1 14 1 3 24 6 6 9 5 2 8 -1 5 3 3 -3 1 2 4 2 3 -1 4 10 1 11 15 0 -5 9 6 3 2 -2 0 3 2 1 2 33 2 1 3 4 -2 4 58 2 10 2 3 5 3 3 10 3 2 18 3 1 2 -1 7 2 1 1 14 3 7 2 8 5 3 6 2 3 5 2 -1 6 6 3 8 0 3 12 30 2 3 2 6 0 17'
$wsYDS.Range("B3").Value = '14 3 17 2 3 3 8 5 6 4 10 8 55 9 7 4 16 3 7 4 27 3 8 23 9 12 17 11 11 6 13 4 9 7 10 12 7 -3 11 9 6 18 2 11 6 5 6 19 15 58 9 21 15 3 16 8 4 11 72 20 11 10 12 49 14 -4 11 12 6 32 19 13 8 6 23 7 14 5 9 14 6 31 15 -5 6 7 13 8 11 13 6 4 5 18 7 16 5 14 12 14 -1 15 1 38 2 14 5 9 4 16 4 -1 9 2 19 5 14 7 19 11 20 4 15 -1 14 19 9 7 6 5 9 25 5 7 22 8 8 10 8 6 7 6 46 11 3 12 21 21 8 10 34 26 22 5 17 15 5 7 6 18 19 5 10 9 26 21 -6 17 2 5 2 10 1 20 6 5 22 4 7 8 14 13 9 7 10 17 9 24 11 -3 7 6 3 7 6 12 -2 -6 25 27 9 24 9 1 9 30 6 1 5 13 13 70 19 4 2 15 13 6 8 6 18 19 6 21 7 26 5 8 9 7 40 6 7 33 8 4 3 17 5 7'
$wsYDS.Range("C3").Value = '4 6 7 14 4 31 10 19 7 8 5 6 2 7 32 20 20 19 8 5 13 18 19 5 16 7 17 2 10 6 7 13 8 1 4 8 11 11 16 9 -1 14 7 2 9 6 17 31 13 7 3 10 8 27 7 9 4 4 4 4 11 5 22 9 10 4 5 -5 10 13 -4 6 18 4 27 51 19 7 8 17 20 8 4 5 21 1 54 8 5 11 3 28 17 8 18 1 23 20 9 18 8 7 17 13 7 16 10 7 -7 9 3 6 4 6 9 26 27 84 5 -4 5 10 10 1 7 12 9 28 0 8 25 12 7 8 8 4 20 31 4 9 7 44 6 7 16 6 6 41 9 7 -1 8 2 12 8 50 11 34 7 39 13 6 0 8 7 5 16 7 20 7 6 3 12 49 16 13 -4 9 64 3 2 50 9 16 17 10 1 16 16 9 10 -1 2 5 7 2 39 0 8 4 16 14 33 -2 4 14 6 14 14 -1 23 14 9 7 7 26 1 6 12 7 9 4 17 4 7 14 5 21 4 5 6 28 23 9 5 1 41 7 1 21 10 10 5 17 1 4 20 0 12 2 10 1 3 4 13 25 9 4 3 8 3 2 22 16 11 9 1 6 6 27 3 17 14 14 8 13 5 7 7 -2 3 23 5 8 1 5 5 6 14 19 7 5 8 5 16 6 10 3 6 18 11 13'

# --- OFF sheet: updated RATT (row2) and PATT (row3) totals ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 180
$wsOFF.Range("E2").Value = 11
$wsOFF.Range("F2").Value = 61
$wsOFF.Range("G2").Value = 67
$wsOFF.Range("H2").Value = 8
$wsOFF.Range("I2").Value = 13
$wsOFF.Range("J2").Value = 26
$wsOFF.Range("L2").Value = 226
$wsOFF.Range("M2").Value = 130
$wsOFF.Range("Q2").Value = 444
$wsOFF.Range("B3").Value = 12
$wsOFF.Range("C3").Value = 140
$wsOFF.Range("E3").Value = 34
$wsOFF.Range("F3").Value = 85
$wsOFF.Range("G3").Value = 17
$wsOFF.Range("H3").Value = 25
$wsOFF.Range("I3").Value = 51
$wsOFF.Range("J3").Value = 40
$wsOFF.Range("N3").Value = 8

# --- DEF sheet: updated RATT (row2) and PATT (row3) totals ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 170
$wsDEF.Range("F2").Value = 51
$wsDEF.Range("G2").Value = 40
$wsDEF.Range("I2").Value = 3
$wsDEF.Range("J2").Value = 27
$wsDEF.Range("L2").Value = 254
$wsDEF.Range("M2").Value = 168
$wsDEF.Range("O2").Value = 17
$wsDEF.Range("Q2").Value = 410
$wsDEF.Range("C3").Value = 151
$wsDEF.Range("E3").Value = 30
$wsDEF.Range("F3").Value = 104
$wsDEF.Range("G3").Value = 26
$wsDEF.Range("I3").Value = 66
$wsDEF.Range("J3").Value = 55
$wsDEF.Range("N3").Value = 13

# --- ST sheet: kicking/special-teams totals + per-game strings ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 66
$wsST.Range("D2").Value = 68
$wsST.Range("J2").Value = 42
$wsST.Range("K2").Value = 38
$wsST.Range("L2").Value = 18
$wsST.Range("M2").Value = 10
$wsST.Range("B3").Value = 41
$wsST.Range("B4").Value = '57 65 64 66 56 65 60 63 65 63 65 61 62 47 64 63 59 68 69 65 68 59 62 65 60'
$wsST.Range("B5").Value = '17 15 22 26 17 23 13 17 19 21 25 27 20 7 18 14 20 22 30 23 21 23 42 17 22'
$wsST.Range("D3").Value = '45 58 44 60 42 52 59 37 51 44 54 56 44 34 55 43 46 60 57 52 53 58 41 41 48 41 56 45 50 63 39 45 38 39 38 41 41 43 62 53 39 54 49 42 52 42 43 53 47 59 52 55 33 62 45 62 48 49 54 45 58 33 46 54 48 46 49 33'
$wsST.Range("D4").Value = '0 17 0 10 14 0 0 0 0 8 0 18 12 0 26 0 0 0 0 10 0 8 0 0 0 0 11 0 3 11 0 0 0 0 0 9 0 0 12 10 0 8 0 13 11 10 1 18 7 16 9 0 0 12 0 28 5 0 11 6 11 3 2 9 10 0 0 0'
$wsST.Range("D5").Value = '0 0 13 0 0 0 0 0 25 0 0 0 0 9 0 0 0 13 14 2 0 2 0 0 0 0 0 12 0 0 17 7 26 0 4 0 0 0 8 0 14 7 12 0 8 19 0 8 0 2 0 11 12 13 0 0 18 0 34 0 0 4 8 0 0'

# --- TURNS sheet: turnover counts ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("C2").Value = 9
$wsTURNS.Range("D2").Value = 3
$wsTURNS.Range("E2").Value = 5
$wsTURNS.Range("D3").Value = 2

# --- PEN sheet: penalty counts ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 15
$wsPEN.Range("B3").Value = 16
$wsPEN.Range("D4").Value = 10

